$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Parent company" (B) and "Location County/City" (old E) columns.
# Deleting B first shifts "Location County/City" from E to D, so delete D next.
$ws.Range("B:B").EntireColumn.Delete()
$ws.Range("D:D").EntireColumn.Delete()

# Match the author's final selection on the data row.
$ws.Range("A2:I2").Select()
